# Generate Report for handoff
# b.md has been re-handed-off for both zh-cn and de-de:
#   - Status changes from "Handed back: in sync with en-US" to "Ready for handoff"
#   - Latest Handoff File / Latest Handoff Datetime are updated to the new
#     handoff package for b.md
# This touches the Overview sheet (summary) plus the per-locale zh-cn / de-de
# sheets (detail rows), row 3 in each case (the b.md row).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: collapse status for b.md row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: b.md row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-25 03:28:26"

foreach ($h in $zhcn.Hyperlinks) {
    if ($h.Range.Row -eq 3 -and $h.Range.Column -eq 3) {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# --- de-de sheet: b.md row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-25 03:28:36"

foreach ($h in $dede.Hyperlinks) {
    if ($h.Range.Row -eq 3 -and $h.Range.Column -eq 3) {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
